# Remove the "Reference" slide (slide 1) and the "Customize this Template"
# slide (the last slide), keeping the "DOT NET CORE RABBIT MQ PRODUCER"
# slide and the screenshot slide that follows it, in their original order.

$p = $ppt.ActivePresentation

# Delete from the end first so earlier indices remain stable while we work.
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()

# Remove the original first slide ("Reference").
$p.Slides.Item(1).Delete()
